$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New table: LEACH vs O_BLEACH results (rows 14-19, columns B-H) ---
# The order in which cells are written below matters: it reproduces the
# original author's shared-string insertion order in sharedStrings.xml.

# Header row (row 14) - first three header cells
$ws.Range("B14").Value = "spread"
$ws.Range("C14").Value = "LEACH"
$ws.Range("D14").Value = "O_BLEACH"

# "spread" values going down column B for the data rows.
# These are stored as text (same as the existing 0.5/0.4/0.3/0.2 cells
# B6:B9), so copy the existing text cells instead of typing numbers.
$ws.Range("B6").Copy()
$ws.Range("B15").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("B7").Copy()
$ws.Range("B16").PasteSpecial(-4163)
$ws.Range("B8").Copy()
$ws.Range("B17").PasteSpecial(-4163)
$ws.Range("B9").Copy()
$ws.Range("B18").PasteSpecial(-4163)

# B19 = "0.1" is new text that does not exist anywhere yet. Build it with
# a helper formula and paste its value so it lands as text, not a number.
$ws.Range("Z1").Formula = '="0.1"'
$ws.Range("Z1").Copy()
$ws.Range("B19").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

# LEACH / O_BLEACH packet counts (columns C & D)
$ws.Range("C15").Value = 5742
$ws.Range("D15").Value = 5767
$ws.Range("C16").Value = 6117
$ws.Range("D16").Value = 6132
$ws.Range("C17").Value = 5961
$ws.Range("D17").Value = 5973
$ws.Range("C18").Value = 5920
$ws.Range("D18").Value = 5898
$ws.Range("C19").Value = 5706
$ws.Range("D19").Value = 5655

# Finish header row (row 14)
$ws.Range("E14").Value = "dataSent LEACH [packets]"
$ws.Range("F14").Value = "dataSent BLEACH [packets]"
$ws.Range("G14").Value = "startNRJ"
$ws.Range("H14").Value = "Runs"

# Row 15: dataSent / startNRJ / Runs
$ws.Range("F15").Value = "159, 316"
$ws.Range("E15").Value = "159, 390"
$ws.Range("B6").Copy()
$ws.Range("G15").PasteSpecial(-4163)
$ws.Range("H15").Value = 50

# Row 16
$ws.Range("E16").Value = "157, 238"
$ws.Range("F16").Value = "157, 244"
$ws.Range("B6").Copy()
$ws.Range("G16").PasteSpecial(-4163)
$ws.Range("H16").Value = 50

# Row 17
$ws.Range("E17").Value = "155, 648"
$ws.Range("F17").Value = "155, 700"
$ws.Range("B6").Copy()
$ws.Range("G17").PasteSpecial(-4163)
$ws.Range("H17").Value = 50

# Row 18
$ws.Range("F18").Value = "158, 451"
$ws.Range("E18").Value = "158, 479"
$ws.Range("B6").Copy()
$ws.Range("G18").PasteSpecial(-4163)
$ws.Range("H18").Value = 50

# Row 19
$ws.Range("F19").Value = "156, 711"
$ws.Range("E19").Value = "156, 591"
$ws.Range("B6").Copy()
$ws.Range("G19").PasteSpecial(-4163)
$ws.Range("H19").Value = 50

# Right-align the Runs address columns (E15:F19), matching the style
# already used by E3:F3 in the existing table.
$ws.Range("E15:F19").HorizontalAlignment = -4152

# Update the view: select the last-entered cell.
$ws.Range("E19").Select()
